# This script reproduces the "solution generator" edit:
#  1) On sheets "Miesiac 3", "Miesiac 7", "Miesiac 11" a handful of result
#     cells (columns B/C/D) get new (non-zero) values.
#  2) On every monthly sheet ("Miesiac 1".."Miesiac 12") the two grape
#     columns "Arneis" (E) and "Cortese" (F) are removed entirely - the
#     cells are deleted and everything to the right shifts left, shrinking
#     each sheet's used range from A1:F4 down to A1:D4.
#  3) Once no sheet references the "Arneis"/"Cortese" shared strings any
#     more, they naturally drop out of the shared string table.

$wb = $excel.ActiveWorkbook

# --- Step 1: update the result values on the affected monthly sheets ---

$ws3 = $wb.Worksheets.Item("Miesiac 3")
$ws3.Range("C2").Value = 400
$ws3.Range("B3").Value = 400
$ws3.Range("B4").Value = 400

$ws7 = $wb.Worksheets.Item("Miesiac 7")
$ws7.Range("C2").Value = 400
$ws7.Range("B3").Value = 400
$ws7.Range("C3").Value = 0
$ws7.Range("B4").Value = 400
$ws7.Range("D4").Value = 0

$ws11 = $wb.Worksheets.Item("Miesiac 11")
$ws11.Range("B2").Value = 0
$ws11.Range("C2").Value = 400
$ws11.Range("B3").Value = 400
$ws11.Range("D3").Value = 0
$ws11.Range("B4").Value = 400

# --- Step 2: drop the "Arneis" (E) and "Cortese" (F) columns everywhere ---

foreach ($ws in $wb.Worksheets) {
    $ws.Range("E1:F4").Delete()
}
